$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns G:K (numeric-looking values) to be stored as text,
# matching the source data which keeps these as literal strings (e.g. "75.00", "0").
$ws.Range("G14:K25").NumberFormat = "@"

$rows = @(
    @(" Abu Dhabi", " October 25 2020", "Royals won by 8 wickets (with 10 balls remaining)", "Mumbai Indians", "Rajasthan Royals", "Krunal Pandya ", "3", "4", "0", "0", "75.00"),
    @(" Abu Dhabi", " September 19 2020", "Super Kings won by 5 wickets (with 4 balls remaining)", "Mumbai Indians", "Chennai Super Kings", "Krunal Pandya ", "3", "3", "0", "0", "100.00"),
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Mumbai Indians", "Royal Challengers Bangalore", "Krunal Pandya ", "10", "10", "1", "0", "100.00"),
    @(" Abu Dhabi", " September 23 2020", "Mumbai won by 49 runs", "Mumbai Indians", "Kolkata Knight Riders", "Krunal Pandya ", "1", "3", "0", "0", "33.33"),
    @(" Abu Dhabi", " October 11 2020", "Mumbai won by 5 wickets (with 2 balls remaining)", "Mumbai Indians", "Delhi Capitals", "Krunal Pandya ", "12", "7", "2", "0", "171.42"),
    @(" Dubai (DSC)", " November 05 2020", "Mumbai won by 57 runs", "Mumbai Indians", "Delhi Capitals", "Krunal Pandya ", "13", "10", "0", "1", "130.00"),
    @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Mumbai Indians", "Sunrisers Hyderabad", "Krunal Pandya ", "20", "4", "2", "2", "500.00"),
    @(" Sharjah", " November 03 2020", "Sunrisers won by 10 wickets (with 17 balls remaining)", "Mumbai Indians", "Sunrisers Hyderabad", "Krunal Pandya ", "0", "3", "0", "0", "0.00"),
    @(" Dubai (DSC)", " November 10 2020", "Mumbai won by 5 wickets (with 8 balls remaining)", "Mumbai Indians", "Delhi Capitals", "Krunal Pandya ", "1", "1", "0", "0", "100.00"),
    @(" Dubai (DSC)", " September 28 2020", "Match tied (RCB won the one-over eliminator)", "Mumbai Indians", "Royal Challengers Bangalore", "Krunal Pandya ", "0", "0", "0", "0", "-"),
    @(" Abu Dhabi", " October 06 2020", "Mumbai won by 57 runs", "Mumbai Indians", "Rajasthan Royals", "Krunal Pandya ", "12", "17", "0", "1", "70.58"),
    @(" Dubai (DSC)", " October 18 2020", "Match tied (Kings XI won the one-over eliminator)", "Mumbai Indians", "Kings XI Punjab", "Krunal Pandya ", "34", "30", "4", "1", "113.33"),
)

$r = 14
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
